$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Title
$meta.Range("B5").Value = "DMI Code EMDN"

# Date
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"

# Description
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter le code EMDN."

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Short
$elements.Range("L2").Value = "DMI Code EMDN"

# Definition
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter le code EMDN."

# Mapping: RIM Mapping (cleared)
$elements.Range("AK2").Value = ""
